# Adds a new "2021" column (R) of data to the table, mirroring the
# formatting of the existing "2020" column (Q), and moves the active
# selection to O10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R-column values for rows 4-44 (row => new 2021 value). Rows that are
# intentionally blank separators in the source table (6, 9, 12, 17) are
# omitted here and handled by the format-only pass below.
$values = [ordered]@{
  4  = 2021
  5  = 33.262233298138462
  7  = 33.10733359588
  8  = 33.434791912906064
  10 = 33.257744153601877
  11 = 33.264901349775037
  13 = 40.506761643955343
  14 = 31.172874511993772
  15 = 29.650854888195603
  16 = 24.097659040517694
  18 = 40.671537678982844
  19 = 40.835744431591088
  20 = 40.494322790314847
  21 = 43.213901910043809
  22 = 43.25339669708363
  23 = 43.172054336673064
  24 = 38.130315382405762
  25 = 39.052986923894757
  26 = 37.162402419999459
  27 = 39.151777291250745
  28 = 39.257524172776719
  29 = 39.034895273433577
  30 = 23.758244663001044
  31 = 23.600989554960133
  32 = 23.928051635532135
  33 = 23.479873624436866
  34 = 24.191255182099319
  35 = 22.759796314587014
  36 = 27.00352039684709
  37 = 26.49991066711625
  38 = 27.572066513890491
  39 = 35.769118603355516
  40 = 35.037028356467729
  41 = 36.714712776303102
  42 = 28.585005142864613
  43 = 27.745516596290607
  44 = 29.549901683892372
}

# Rows where R gets the column's formatting but stays empty (blank
# separator rows in the source data, same as column Q there).
$blankRows = @(6, 9, 12, 17)

$xlPasteFormats = -4122

for ($row = 4; $row -le 44; $row++) {
    $ws.Range("Q$row").Copy() | Out-Null
    $ws.Range("R$row").PasteSpecial($xlPasteFormats) | Out-Null
}
$excel.CutCopyMode = $false

foreach ($row in $values.Keys) {
    $ws.Range("R$row").Value = $values[$row]
}

# Move the active selection to O10, as recorded in the sheet view.
$ws.Range("O10").Select() | Out-Null
